# Update the 100 arithmetic answers in the "within 100" worksheet table.
# Each old equation string is unique in the document, so a scoped
# Find/Replace (MatchWholeWord, Replace=2/wdReplaceAll) per pair is safe.
$d = $word.ActiveDocument
$d.Content.Find.Execute("18+9=27", $true, $true, $false, $false, $false, $true, 1, $false, "12+52=64", 2) | Out-Null
$d.Content.Find.Execute("46-33=13", $true, $true, $false, $false, $false, $true, 1, $false, "28+46=74", 2) | Out-Null
$d.Content.Find.Execute("87-87=0", $true, $true, $false, $false, $false, $true, 1, $false, "44-7=37", 2) | Out-Null
$d.Content.Find.Execute("18-11=7", $true, $true, $false, $false, $false, $true, 1, $false, "27+29=56", 2) | Out-Null
$d.Content.Find.Execute("49+16=65", $true, $true, $false, $false, $false, $true, 1, $false, "48-34=14", 2) | Out-Null
$d.Content.Find.Execute("43+54=97", $true, $true, $false, $false, $false, $true, 1, $false, "22-11=11", 2) | Out-Null
$d.Content.Find.Execute("76-52=24", $true, $true, $false, $false, $false, $true, 1, $false, "32+22=54", 2) | Out-Null
$d.Content.Find.Execute("91-48=43", $true, $true, $false, $false, $false, $true, 1, $false, "32+62=94", 2) | Out-Null
$d.Content.Find.Execute("35-25=10", $true, $true, $false, $false, $false, $true, 1, $false, "89-12=77", 2) | Out-Null
$d.Content.Find.Execute("69-36=33", $true, $true, $false, $false, $false, $true, 1, $false, "60-59=1", 2) | Out-Null
$d.Content.Find.Execute("27+34=61", $true, $true, $false, $false, $false, $true, 1, $false, "38+40=78", 2) | Out-Null
$d.Content.Find.Execute("46+53=99", $true, $true, $false, $false, $false, $true, 1, $false, "66-54=12", 2) | Out-Null
$d.Content.Find.Execute("22+15=37", $true, $true, $false, $false, $false, $true, 1, $false, "59+2=61", 2) | Out-Null
$d.Content.Find.Execute("6+41=47", $true, $true, $false, $false, $false, $true, 1, $false, "30+10=40", 2) | Out-Null
$d.Content.Find.Execute("58-48=10", $true, $true, $false, $false, $false, $true, 1, $false, "53-1=52", 2) | Out-Null
$d.Content.Find.Execute("45-10=35", $true, $true, $false, $false, $false, $true, 1, $false, "88-11=77", 2) | Out-Null
$d.Content.Find.Execute("36+2=38", $true, $true, $false, $false, $false, $true, 1, $false, "45+0=45", 2) | Out-Null
$d.Content.Find.Execute("16+81=97", $true, $true, $false, $false, $false, $true, 1, $false, "64-1=63", 2) | Out-Null
$d.Content.Find.Execute("20+16=36", $true, $true, $false, $false, $false, $true, 1, $false, "58-10=48", 2) | Out-Null
$d.Content.Find.Execute("76-45=31", $true, $true, $false, $false, $false, $true, 1, $false, "18+7=25", 2) | Out-Null
$d.Content.Find.Execute("64-14=50", $true, $true, $false, $false, $false, $true, 1, $false, "9+60=69", 2) | Out-Null
$d.Content.Find.Execute("8+38=46", $true, $true, $false, $false, $false, $true, 1, $false, "18+76=94", 2) | Out-Null
$d.Content.Find.Execute("11+22=33", $true, $true, $false, $false, $false, $true, 1, $false, "43+49=92", 2) | Out-Null
$d.Content.Find.Execute("77-50=27", $true, $true, $false, $false, $false, $true, 1, $false, "12+4=16", 2) | Out-Null
$d.Content.Find.Execute("2+24=26", $true, $true, $false, $false, $false, $true, 1, $false, "66+1=67", 2) | Out-Null
$d.Content.Find.Execute("50+48=98", $true, $true, $false, $false, $false, $true, 1, $false, "19+60=79", 2) | Out-Null
$d.Content.Find.Execute("96-80=16", $true, $true, $false, $false, $false, $true, 1, $false, "23-12=11", 2) | Out-Null
$d.Content.Find.Execute("95-88=7", $true, $true, $false, $false, $false, $true, 1, $false, "6+74=80", 2) | Out-Null
$d.Content.Find.Execute("9+39=48", $true, $true, $false, $false, $false, $true, 1, $false, "69-34=35", 2) | Out-Null
$d.Content.Find.Execute("22-2=20", $true, $true, $false, $false, $false, $true, 1, $false, "63-63=0", 2) | Out-Null
$d.Content.Find.Execute("68-24=44", $true, $true, $false, $false, $false, $true, 1, $false, "87-57=30", 2) | Out-Null
$d.Content.Find.Execute("40+4=44", $true, $true, $false, $false, $false, $true, 1, $false, "27+23=50", 2) | Out-Null
$d.Content.Find.Execute("95-46=49", $true, $true, $false, $false, $false, $true, 1, $false, "69-33=36", 2) | Out-Null
$d.Content.Find.Execute("32+43=75", $true, $true, $false, $false, $false, $true, 1, $false, "49-27=22", 2) | Out-Null
$d.Content.Find.Execute("19+8=27", $true, $true, $false, $false, $false, $true, 1, $false, "12-9=3", 2) | Out-Null
$d.Content.Find.Execute("5+85=90", $true, $true, $false, $false, $false, $true, 1, $false, "54+25=79", 2) | Out-Null
$d.Content.Find.Execute("73-0=73", $true, $true, $false, $false, $false, $true, 1, $false, "50-0=50", 2) | Out-Null
$d.Content.Find.Execute("68+11=79", $true, $true, $false, $false, $false, $true, 1, $false, "23+16=39", 2) | Out-Null
$d.Content.Find.Execute("95-74=21", $true, $true, $false, $false, $false, $true, 1, $false, "93-85=8", 2) | Out-Null
$d.Content.Find.Execute("14+13=27", $true, $true, $false, $false, $false, $true, 1, $false, "27-3=24", 2) | Out-Null
$d.Content.Find.Execute("56-19=37", $true, $true, $false, $false, $false, $true, 1, $false, "39+18=57", 2) | Out-Null
$d.Content.Find.Execute("89+9=98", $true, $true, $false, $false, $false, $true, 1, $false, "28-18=10", 2) | Out-Null
$d.Content.Find.Execute("34+51=85", $true, $true, $false, $false, $false, $true, 1, $false, "7+55=62", 2) | Out-Null
$d.Content.Find.Execute("95-11=84", $true, $true, $false, $false, $false, $true, 1, $false, "42+19=61", 2) | Out-Null
$d.Content.Find.Execute("25+27=52", $true, $true, $false, $false, $false, $true, 1, $false, "84-66=18", 2) | Out-Null
$d.Content.Find.Execute("83-47=36", $true, $true, $false, $false, $false, $true, 1, $false, "59-37=22", 2) | Out-Null
$d.Content.Find.Execute("27-20=7", $true, $true, $false, $false, $false, $true, 1, $false, "83-39=44", 2) | Out-Null
$d.Content.Find.Execute("0+54=54", $true, $true, $false, $false, $false, $true, 1, $false, "98-9=89", 2) | Out-Null
$d.Content.Find.Execute("51-49=2", $true, $true, $false, $false, $false, $true, 1, $false, "75-17=58", 2) | Out-Null
$d.Content.Find.Execute("81+3=84", $true, $true, $false, $false, $false, $true, 1, $false, "15-8=7", 2) | Out-Null
$d.Content.Find.Execute("90-5=85", $true, $true, $false, $false, $false, $true, 1, $false, "43-15=28", 2) | Out-Null
$d.Content.Find.Execute("20+42=62", $true, $true, $false, $false, $false, $true, 1, $false, "16+41=57", 2) | Out-Null
$d.Content.Find.Execute("87-18=69", $true, $true, $false, $false, $false, $true, 1, $false, "47-18=29", 2) | Out-Null
$d.Content.Find.Execute("69-44=25", $true, $true, $false, $false, $false, $true, 1, $false, "47+9=56", 2) | Out-Null
$d.Content.Find.Execute("76-64=12", $true, $true, $false, $false, $false, $true, 1, $false, "41+18=59", 2) | Out-Null
$d.Content.Find.Execute("85+5=90", $true, $true, $false, $false, $false, $true, 1, $false, "3+47=50", 2) | Out-Null
$d.Content.Find.Execute("25-24=1", $true, $true, $false, $false, $false, $true, 1, $false, "70-64=6", 2) | Out-Null
$d.Content.Find.Execute("62-20=42", $true, $true, $false, $false, $false, $true, 1, $false, "12-7=5", 2) | Out-Null
$d.Content.Find.Execute("96-50=46", $true, $true, $false, $false, $false, $true, 1, $false, "15+74=89", 2) | Out-Null
$d.Content.Find.Execute("97-28=69", $true, $true, $false, $false, $false, $true, 1, $false, "41+52=93", 2) | Out-Null
$d.Content.Find.Execute("27-26=1", $true, $true, $false, $false, $false, $true, 1, $false, "36+58=94", 2) | Out-Null
$d.Content.Find.Execute("10+69=79", $true, $true, $false, $false, $false, $true, 1, $false, "35-32=3", 2) | Out-Null
$d.Content.Find.Execute("31+67=98", $true, $true, $false, $false, $false, $true, 1, $false, "31+37=68", 2) | Out-Null
$d.Content.Find.Execute("1+13=14", $true, $true, $false, $false, $false, $true, 1, $false, "92-44=48", 2) | Out-Null
$d.Content.Find.Execute("79-19=60", $true, $true, $false, $false, $false, $true, 1, $false, "89-35=54", 2) | Out-Null
$d.Content.Find.Execute("80-79=1", $true, $true, $false, $false, $false, $true, 1, $false, "36-13=23", 2) | Out-Null
$d.Content.Find.Execute("2+13=15", $true, $true, $false, $false, $false, $true, 1, $false, "58-1=57", 2) | Out-Null
$d.Content.Find.Execute("3+54=57", $true, $true, $false, $false, $false, $true, 1, $false, "45+13=58", 2) | Out-Null
$d.Content.Find.Execute("9+10=19", $true, $true, $false, $false, $false, $true, 1, $false, "7+35=42", 2) | Out-Null
$d.Content.Find.Execute("2+68=70", $true, $true, $false, $false, $false, $true, 1, $false, "93-9=84", 2) | Out-Null
$d.Content.Find.Execute("84-37=47", $true, $true, $false, $false, $false, $true, 1, $false, "64+17=81", 2) | Out-Null
$d.Content.Find.Execute("59-21=38", $true, $true, $false, $false, $false, $true, 1, $false, "65-16=49", 2) | Out-Null
$d.Content.Find.Execute("52-9=43", $true, $true, $false, $false, $false, $true, 1, $false, "47+23=70", 2) | Out-Null
$d.Content.Find.Execute("10+88=98", $true, $true, $false, $false, $false, $true, 1, $false, "73-69=4", 2) | Out-Null
$d.Content.Find.Execute("54-14=40", $true, $true, $false, $false, $false, $true, 1, $false, "73+15=88", 2) | Out-Null
$d.Content.Find.Execute("48-40=8", $true, $true, $false, $false, $false, $true, 1, $false, "45-21=24", 2) | Out-Null
$d.Content.Find.Execute("78+11=89", $true, $true, $false, $false, $false, $true, 1, $false, "15+38=53", 2) | Out-Null
$d.Content.Find.Execute("1+74=75", $true, $true, $false, $false, $false, $true, 1, $false, "40+48=88", 2) | Out-Null
$d.Content.Find.Execute("19-15=4", $true, $true, $false, $false, $false, $true, 1, $false, "98-71=27", 2) | Out-Null
$d.Content.Find.Execute("97-47=50", $true, $true, $false, $false, $false, $true, 1, $false, "65-43=22", 2) | Out-Null
$d.Content.Find.Execute("41+54=95", $true, $true, $false, $false, $false, $true, 1, $false, "37+43=80", 2) | Out-Null
$d.Content.Find.Execute("30-23=7", $true, $true, $false, $false, $false, $true, 1, $false, "31+17=48", 2) | Out-Null
$d.Content.Find.Execute("97-68=29", $true, $true, $false, $false, $false, $true, 1, $false, "32+66=98", 2) | Out-Null
$d.Content.Find.Execute("5+43=48", $true, $true, $false, $false, $false, $true, 1, $false, "8+78=86", 2) | Out-Null
$d.Content.Find.Execute("74-49=25", $true, $true, $false, $false, $false, $true, 1, $false, "89-76=13", 2) | Out-Null
$d.Content.Find.Execute("40+18=58", $true, $true, $false, $false, $false, $true, 1, $false, "70-49=21", 2) | Out-Null
$d.Content.Find.Execute("78-49=29", $true, $true, $false, $false, $false, $true, 1, $false, "92-60=32", 2) | Out-Null
$d.Content.Find.Execute("67+23=90", $true, $true, $false, $false, $false, $true, 1, $false, "33+2=35", 2) | Out-Null
$d.Content.Find.Execute("5+54=59", $true, $true, $false, $false, $false, $true, 1, $false, "7+42=49", 2) | Out-Null
$d.Content.Find.Execute("50+28=78", $true, $true, $false, $false, $false, $true, 1, $false, "30+36=66", 2) | Out-Null
$d.Content.Find.Execute("35+37=72", $true, $true, $false, $false, $false, $true, 1, $false, "89-36=53", 2) | Out-Null
$d.Content.Find.Execute("8+10=18", $true, $true, $false, $false, $false, $true, 1, $false, "12+38=50", 2) | Out-Null
$d.Content.Find.Execute("44-9=35", $true, $true, $false, $false, $false, $true, 1, $false, "11-0=11", 2) | Out-Null
$d.Content.Find.Execute("65+28=93", $true, $true, $false, $false, $false, $true, 1, $false, "6+5=11", 2) | Out-Null
$d.Content.Find.Execute("47-33=14", $true, $true, $false, $false, $false, $true, 1, $false, "87-11=76", 2) | Out-Null
$d.Content.Find.Execute("98-83=15", $true, $true, $false, $false, $false, $true, 1, $false, "55-19=36", 2) | Out-Null
$d.Content.Find.Execute("0+56=56", $true, $true, $false, $false, $false, $true, 1, $false, "95-45=50", 2) | Out-Null
$d.Content.Find.Execute("1+6=7", $true, $true, $false, $false, $false, $true, 1, $false, "76-74=2", 2) | Out-Null
$d.Content.Find.Execute("14+77=91", $true, $true, $false, $false, $false, $true, 1, $false, "25+30=55", 2) | Out-Null
$d.Content.Find.Execute("49+14=63", $true, $true, $false, $false, $false, $true, 1, $false, "41+29=70", 2) | Out-Null
